$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.398.60'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.238.86'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  -0.07%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '243.74'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('E6').Value = '  -0.35%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '74.50'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -2.60%  '
$ws.Range('E8').Value = '  +0.08%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.612'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -2.51%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '42.85'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -4.15%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.0967'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('E12').Value = '  -4.56%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = '2.576.60'
$ws.Range('E14').Value = '  -0.35%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '14.40'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -2.43%  '
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '2.230.66'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '42.245.79'
$ws.Range('E18').Value = '  -0.48%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '0.0000106'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +4.59%  '
$ws.Range('E20').Value = '  +0.58%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '73.26'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +1.37%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '11.26'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  +1.37%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '231.79'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -0.26%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '2.11'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -6.83%  '
$ws.Range('E25').Value = '  +0.10%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '11.47'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('E29').Value = '  -2.06%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '167.29'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -0.08%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '20.65'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -0.41%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '5.74'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('E33').Value = '  -2.35%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '30.44'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -6.63%  '
$ws.Range('E35').Value = '  -0.63%  '
$ws.Range('E36').Value = '  -9.32%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '4.38'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -8.13%  '
$ws.Range('E38').Value = '  -4.08%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '13.61'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -4.64%  '
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('E41').Value = '  -1.66%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '65.06'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('E43').Value = '  -1.98%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '8.76'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -1.97%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '105.25'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -2.79%  '
$ws.Range('E46').Value = '  -2.10%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '2.37'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -1.22%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.13'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -1.86%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.18'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').Value = '2.447.47'
$ws.Range('E51').Value = '  -0.80%  '
